# Applies the "added parameter and data summary" edit:
# - results sheet gains a "data" block (A:D) beside the existing "parameter" block (E:H)
# - plot sheet column D (predicted stress) and loss sheet column B (loss history) get new values

$wb = $excel.ActiveWorkbook
$results = $wb.Worksheets.Item("results")
$plot = $wb.Worksheets.Item("plot")
$loss = $wb.Worksheets.Item("loss")

# ---------------------------------------------------------------
# results sheet: rebuild A1:H6 as two side-by-side tables
#   A:D -> new "data" table (time/strain/stress/temperature/cycle)
#   E:H -> old "parameter" table (n/eta/s0/R/d), now shifted right
# ---------------------------------------------------------------

# header row (bold, bordered, centered - same style as before)
$results.Range("A1").Value = "data"
$results.Range("B1").Value = "d. bounds"
$results.Range("C1").Value = "d. scales"
$results.Range("D1").Value = "|"
$results.Range("E1").Value = "parameter"
$results.Range("F1").Value = "p. bounds"
$results.Range("G1").Value = "p. scales"
$results.Range("H1").Value = "p. optimised"
$results.Range("A1:H1").Font.Bold = $true
$results.Range("A1:H1").Borders.LineStyle = 1
$results.Range("A1:H1").HorizontalAlignment = -4108
$results.Range("A1:H1").VerticalAlignment = -4160

# data rows (centered, unbordered - same style as before)
$results.Range("A2").Value = "time"
$results.Range("B2").Value = "[0.0, 50.0]"
$results.Range("C2").Value = "[0.0, 1.0]"
$results.Range("D2").Value = "|"
$results.Range("E2").Value = "n"
$results.Range("F2").Value = "[3.5, 10.5]"
$results.Range("G2").Value = "[0.0, 10.0]"
$results.Range("H2").Value = 6.660883283623221
$results.Range("A3").Value = "strain"
$results.Range("B3").Value = "[0.0, 0.5]"
$results.Range("C3").Value = "[0.0, 1.0]"
$results.Range("D3").Value = "|"
$results.Range("E3").Value = "eta"
$results.Range("F3").Value = "[1.5e+02, 4.5e+02]"
$results.Range("G3").Value = "[0.0, 10.0]"
$results.Range("H3").Value = 363.4621364717364
$results.Range("A4").Value = "stress"
$results.Range("B4").Value = "[0.0, 3.9e+02]"
$results.Range("C4").Value = "[0.0, 1.0]"
$results.Range("D4").Value = "|"
$results.Range("E4").Value = "s0"
$results.Range("F4").Value = "[25.0, 75.0]"
$results.Range("G4").Value = "[0.0, 1.0]"
$results.Range("H4").Value = -158.2609353223543
$results.Range("A5").Value = "temperature"
$results.Range("B5").Value = "[0.0, 0.0]"
$results.Range("C5").Value = "[0.0, 1.0]"
$results.Range("D5").Value = "|"
$results.Range("E5").Value = "R"
$results.Range("F5").Value = "[1e+02, 3e+02]"
$results.Range("G5").Value = "[0.0, 1.0]"
$results.Range("H5").Value = 202.3720300123277
$results.Range("A6").Value = "cycle"
$results.Range("B6").Value = "[0.0, 0.0]"
$results.Range("C6").Value = "[0.0, 1.0]"
$results.Range("D6").Value = "|"
$results.Range("E6").Value = "d"
$results.Range("F6").Value = "[2.5, 7.5]"
$results.Range("G6").Value = "[0.0, 1.0]"
$results.Range("H6").Value = 2.507102890713211
$results.Range("A2:H6").HorizontalAlignment = -4108

# column widths (A,C widened; new B..H narrow divider + new "data" cols)
$results.Columns.Item(1).ColumnWidth = 12.7109375
$results.Columns.Item(3).ColumnWidth = 11.7109375
$results.Columns.Item(4).ColumnWidth = 2.7109375
$results.Columns.Item(5).ColumnWidth = 10.7109375
$results.Columns.Item(6).ColumnWidth = 19.7109375
$results.Columns.Item(7).ColumnWidth = 12.7109375
$results.Columns.Item(8).ColumnWidth = 20.7109375

# ---------------------------------------------------------------
# plot sheet: refreshed "predicted" series (column D) from the new fit
# ---------------------------------------------------------------
$plot.Range("D3").Value = 210.3386149766334
$plot.Range("D4").Value = 223.2653294220156
$plot.Range("D5").Value = 232.1451983464777
$plot.Range("D6").Value = 240.5403686328063
$plot.Range("D7").Value = 248.5337687489057
$plot.Range("D8").Value = 256.1455444602247
$plot.Range("D9").Value = 263.3938476631907
$plot.Range("D10").Value = 270.2959571812147
$plot.Range("D11").Value = 276.8683373779972
$plot.Range("D12").Value = 283.1266761539799
$plot.Range("D13").Value = 289.0859209530142
$plot.Range("D14").Value = 294.7603131729384
$plot.Range("D15").Value = 300.1634210238941
$plot.Range("D16").Value = 305.3081708936639
$plot.Range("D17").Value = 310.2068772894611
$plot.Range("D18").Value = 314.8712714214943
$plot.Range("D19").Value = 319.3125284877154
$plot.Range("D20").Value = 323.5412937153291
$plot.Range("D21").Value = 327.5677072122668
$plot.Range("D22").Value = 331.4014276799311
$plot.Range("D23").Value = 335.0516550367346
$plot.Range("D24").Value = 338.5271520001228
$plot.Range("D25").Value = 341.836264672967
$plot.Range("D26").Value = 344.9869421784006
$plot.Range("D27").Value = 347.9867553854161
$plot.Range("D28").Value = 350.8429147658038
$plot.Range("D29").Value = 353.5622874213685
$plot.Range("D30").Value = 356.1514133187258
$plot.Range("D31").Value = 358.6165207674283
$plot.Range("D32").Value = 360.9635411756665
$plot.Range("D33").Value = 363.1981231163344
$plot.Range("D34").Value = 365.3256457348441
$plot.Range("D35").Value = 367.3512315287368
$plot.Range("D36").Value = 369.2797585278241
$plot.Range("D37").Value = 371.1158719023549
$plot.Range("D38").Value = 372.8639950254966
$plot.Range("D39").Value = 374.5283400152685
$plot.Range("D40").Value = 376.1129177799445
$plot.Range("D41").Value = 377.621547589904
$plot.Range("D42").Value = 379.0578661978513
$plot.Range("D43").Value = 380.4253362377946
$plot.Range("D44").Value = 381.7272562758239
$plot.Range("D45").Value = 382.966765153736
$plot.Range("D46").Value = 384.1468524219983
$plot.Range("D47").Value = 385.2703646568107
$plot.Range("D48").Value = 386.3400122642378
$plot.Range("D49").Value = 387.3583759757029
$plot.Range("D50").Value = 388.3279130355809
$plot.Range("D51").Value = 389.8871769858835

# ---------------------------------------------------------------
# loss sheet: refreshed loss-history series (column B)
# ---------------------------------------------------------------
$loss.Range("B2").Value = 51.0
$loss.Range("B3").Value = 0.0025
$loss.Range("B4").Value = 0.000068
$loss.Range("B5").Value = 0.000068
$loss.Range("B6").Value = 0.000068
$loss.Range("B7").Value = 0.000068
$loss.Range("B8").Value = 0.000068
$loss.Range("B9").Value = 0.000068
$loss.Range("B10").Value = 0.000068
$loss.Range("B11").Value = 0.000068
